$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2023-09-10 -> 2023-09-11, i.e. Excel serial 45179 -> 45180) for every data
# row from row 2 through row 230.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 230) { $lastRow = 230 }

$ws.Range("C2:C$lastRow").Value = 45180
